$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, pushing the existing rows 42-58 down to 43-59.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly price record
# (same market/category metadata as its neighbours, new date and prices).
$ws.Cells.Item(42, 1).Value = 3
$ws.Cells.Item(42, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44460
$ws.Cells.Item(42, 5).Value = 5
$ws.Cells.Item(42, 6).Value = 100112026
$ws.Cells.Item(42, 7).Value = "Haba"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 40
$ws.Cells.Item(42, 11).Value = 11000
$ws.Cells.Item(42, 12).Value = 11000
$ws.Cells.Item(42, 13).Value = 11000
$ws.Cells.Item(42, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 440
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
